# Fix spelling mistake on the "Twitter OAuth 1.0 authentication" bullet
# (slide 6, content placeholder): merge the 4 split runs
# ("Twitter " / "Oauth" / " 1.0 " / "authentificaiton") into a single,
# correctly spelled run "Twitter OAuth 1.0 authentication".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Locate the last paragraph in the placeholder (the misspelled bullet)
$paraCount = $tr.Paragraphs().Count
$para = $tr.Paragraphs($paraCount, 1)

# Replace its characters (spanning the previously split runs) with the
# corrected, single-run text.
$chars = $tr.Characters($para.Start, $para.Length)
$chars.Text = "Twitter OAuth 1.0 authentication"
